$d = $word.ActiveDocument

# The document contains three "<id>p019v_N</id>" paragraphs that were
# originally split across three runs (one per XML tag / content chunk,
# each with its own formatting). Collapse each one back into a single
# run containing the full "<id>p019v_N</id>" text, picking up the
# surrounding Courier-New / gold-colored formatting that the "<id>" and
# "</id>" runs already use. A literal Find & Replace over the whole
# tag text merges the matched runs into one run automatically.
$ids = @("p019v_1", "p019v_2", "p019v_3")

foreach ($id in $ids) {
    $tag = "<id>" + $id + "</id>"
    $d.Content.Find.Execute($tag, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $tag, 2)
}
